$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44708
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 70
$ws.Range("N2").Value = 12000
$ws.Range("O2").Value = 13000
$ws.Range("P2").Value = 12571
$ws.Range("Q2").Value = "$/caja 12 kilos empedrada"
$ws.Range("R2").Value = "Provincia de Curicó"
$ws.Range("S2").Value = 1048
$ws.Range("T2").Value = 12

# Row 5
$ws.Range("D5").Value = 44742
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 14000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 14500
$ws.Range("Q5").Value = "$/caja 18 kilos granel"
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 806
$ws.Range("T5").Value = 18

# Row 6
$ws.Range("D6").Value = 44714
$ws.Range("N6").Value = 14000
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 14500
$ws.Range("Q6").Value = "$/caja 18 kilos granel"
$ws.Range("S6").Value = 806
$ws.Range("T6").Value = 18

# Row 7
$ws.Range("D7").Value = 44334
$ws.Range("M7").Value = 100
$ws.Range("N7").Value = 11000
$ws.Range("O7").Value = 12000
$ws.Range("P7").Value = 11500
$ws.Range("Q7").Value = "$/caja 12 kilos granel"
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 11500
$ws.Range("T7").Value = 1

# Row 8
$ws.Range("D8").Value = 45084
$ws.Range("N8").Value = 17000
$ws.Range("O8").Value = 18000
$ws.Range("P8").Value = 17500
$ws.Range("R8").Value = "Región del Maule"
$ws.Range("S8").Value = 972

# Row 9
$ws.Range("D9").Value = 44707
$ws.Range("M9").Value = 60
$ws.Range("N9").Value = 12000
$ws.Range("O9").Value = 13000
$ws.Range("P9").Value = 12500
$ws.Range("Q9").Value = "$/caja 12 kilos empedrada"
$ws.Range("R9").Value = "Provincia de Curicó"
$ws.Range("S9").Value = 1042
$ws.Range("T9").Value = 12

# Row 10
$ws.Range("D10").Value = 44330
$ws.Range("N10").Value = 15000
$ws.Range("O10").Value = 16000
$ws.Range("P10").Value = 15500
$ws.Range("R10").Value = "Provincia de Curicó"
$ws.Range("S10").Value = 861

# Row 11 (new row)
$ws.Range("A11").Value = 11
$ws.Range("B11").Value = "Vega Monumental Concepción"
$ws.Range("C11").Value = "Bíobío"
$ws.Range("D11").Value = 44719
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = "Fruta"
$ws.Range("G11").Value = 100107
$ws.Range("H11").Value = "Otros"
$ws.Range("I11").Value = 100107001
$ws.Range("J11").Value = "Caqui"
$ws.Range("K11").Value = "Mankaki"
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 50
$ws.Range("N11").Value = 14000
$ws.Range("O11").Value = 15000
$ws.Range("P11").Value = 14400
$ws.Range("Q11").Value = "$/caja 18 kilos granel"
$ws.Range("R11").Value = "Región del Maule"
$ws.Range("S11").Value = 800
$ws.Range("T11").Value = 18

# Preserve the date number format / style for the new D11 cell, matching the other date cells in column D
$ws.Range("D11").Style = $ws.Range("D10").Style
$ws.Range("D11").NumberFormat = $ws.Range("D10").NumberFormat
